# Regenerate the "K" (strikeouts) column (G) of the save-data sheet.
# The data source for this column changed from a pseudo "Strike#" figure to
# the real strikeout count (K), so every row's value in column G is
# recalculated/rewritten (row 30 happens to keep its original value, 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..58 (header is row 1), in row order.
# (row 30 keeps its original value, 0 -- it was not touched by the change)
$kValues = @(
    1, 0, 2, 0, 0, 2, 2, 2, 0, 1,
    0, 0, 1, 0, 2, 1, 3, 1, 3, 0,
    0, 0, 1, 3, 1, 1, 0, 2, 0, 1,
    1, 0, 1, 0, 1, 2, 0, 1, 0, 1,
    2, 2, 2, 2, 1, 2, 4, 1, 1, 2,
    3, 1, 2, 2, 2, 2, 1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
